$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the E2/E3 "=TRUE()" boolean formulas with the literal text "TRUE",
# while preserving the original cell style (s="1") for both cells.
# A direct Value assignment of "TRUE" is auto-coerced to a Boolean by the
# COM layer (just like real Excel), and a plain quote-prefixed literal
# ("'TRUE") forces text but reassigns a distinct "quote prefix" style - so
# we force the literal text first, then re-apply the original column style
# by copying formats from an untouched, empty cell further down column E
# (which still carries the column's default style).

$ws.Cells.Item(2, 5).Value = "'TRUE"
$ws.Cells.Item(10, 5).Copy() | Out-Null
$ws.Cells.Item(2, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(3, 5).Value = "'TRUE"
$ws.Cells.Item(10, 5).Copy() | Out-Null
$ws.Cells.Item(3, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Range("E2:E3").Select() | Out-Null
